$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "платформа для реализации идеи"
$ws.Range("B3").Value = 45211
$ws.Range("C3").Value = 45212

# Row 4: "создание репозитория"
$ws.Range("B4").Value = 45213
$ws.Range("C4").Value = 45214

# Row 5: "создание проекта "
$ws.Range("B5").Value = 45219
$ws.Range("C5").Value = 45223

# Row 6: "тестирование"
$ws.Range("B6").Value = 45224
$ws.Range("C6").Value = 45225

# Row 7: "решения проблем"
$ws.Range("B7").Value = 45226
$ws.Range("C7").Value = 45228

# Row 8: "сдача проекта"
$ws.Range("B8").Value = 45229
$ws.Range("C8").Value = 45231

# Move / update the active selection to C9, matching the saved view state.
[void]$ws.Range("C9").Select()
